$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "c"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "m"
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "mm"
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "m"
$ws.Range("C4").Value = 3.3
$ws.Rows.Item(4).RowHeight = 15.75

$ws.Range("C1").Select()
